# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 1) "ODI Batting": INNING_NUMBER (col B) was blank for the first two
#    rows -- those cells should not exist at all (no stray empty <c>).
# 2) Add a brand new "ODI Batting Extra" sheet (after "ODI Bowling")
#    holding the freshly scraped per-innings extras.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Part 1: drop the empty INNING_NUMBER cells on "ODI Batting"
# ---------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B3").ClearContents()

# ---------------------------------------------------------------
# Part 2: add "ODI Batting Extra" as the last sheet
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Match this workbook's page margins (same as every other sheet here)
# rather than Excel's out-of-the-box defaults for a freshly added sheet.
$extra.PageSetup.LeftMargin = 54
$extra.PageSetup.RightMargin = 54
$extra.PageSetup.TopMargin = 72
$extra.PageSetup.BottomMargin = 72
$extra.PageSetup.HeaderMargin = 36
$extra.PageSetup.FooterMargin = 36

# Headers
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Match the header styling (bold + border, centered) already used
# throughout this workbook by copying it from an existing header row.
$batting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Data rows. MATCH_CODE values already exist verbatim elsewhere in
# the workbook as plain text cells, so copy them across directly
# (keeps them text-typed, matching the rest of the sheet).
# ---------------------------------------------------------------
$battingExtraSheet = $wb.Worksheets.Item("ODI Batting")

$battingExtraSheet.Range("D2").Copy($extra.Range("A2"))
$battingExtraSheet.Range("D3").Copy($extra.Range("A3"))
$battingExtraSheet.Range("D4").Copy($extra.Range("A4"))

# Numeric batting position
$extra.Range("B2").Value = 11
$extra.Range("B4").Value = 9

# Blank extras (no data scraped yet for these rows)
$battingExtraSheet.Range("B2").Copy($extra.Range("C2"))
$battingExtraSheet.Range("B2").Copy($extra.Range("D2"))
$battingExtraSheet.Range("B2").Copy($extra.Range("E2"))

$battingExtraSheet.Range("B2").Copy($extra.Range("B3"))
$battingExtraSheet.Range("B2").Copy($extra.Range("C3"))
$battingExtraSheet.Range("B2").Copy($extra.Range("D3"))
$battingExtraSheet.Range("B2").Copy($extra.Range("E3"))

# Text values for row 4 + the MAN_OF_MATCH column
$scratch = $extra.Range("Z1")

$scratch.NumberFormat = "@"
$scratch.Value = "0"
$scratch.Copy()
$extra.Range("C4").PasteSpecial(-4163)  # xlPasteValues
$extra.Range("D4").PasteSpecial(-4163)  # xlPasteValues

$scratch.Value = "1.06%"
$scratch.Copy()
$extra.Range("E4").PasteSpecial(-4163)  # xlPasteValues

$scratch.Value = "NO"
$scratch.Copy()
$extra.Range("F2").PasteSpecial(-4163)  # xlPasteValues
$extra.Range("F3").PasteSpecial(-4163)  # xlPasteValues
$extra.Range("F4").PasteSpecial(-4163)  # xlPasteValues

$scratch.ClearContents()
$scratch.NumberFormat = "General"
$scratch.Clear()

$excel.CutCopyMode = $false

# Leave the workbook's original active sheet selected, as it was before
# this edit (adding a sheet otherwise leaves the new one activated).
$wb.Worksheets.Item("Player Info").Activate()

Write-Host "Applied ODI Batting Extra sheet + cleared stray INNING_NUMBER cells"
